# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets.
# F2: 1720 -> 1721
# F3: 29   -> 30
# F4: 481  -> 482
# F5: 158  -> 159
# F7: 650  -> 654

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1721
    $ws.Range("F3").Value = 30
    $ws.Range("F4").Value = 482
    $ws.Range("F5").Value = 159
    $ws.Range("F7").Value = 654
}
